# Updated Fixed-Residue PSSM scores (B2:K21) — supplemental-figure recompute.
# Each string is one row of 10 comma-separated numeric scores for columns B:K.
$rowData = @(
    "-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353",
    "-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353,2.721226090294059,-17.26237017060353,-17.26237017060353",
    "-17.26237017060353,-17.26237017060353,2.478565095340528,-17.26237017060353,3.455025954791026,-17.26237017060353,1.749776933359364,-17.26237017060353,2.043266343322697,-17.26237017060353",
    "-17.26237017060353,4.321919377060516,-17.26237017060353,-17.26237017060353,-17.26237017060353,2.767673205285067,-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353",
    "-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353",
    "2.5427139736264,-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353",
    "-17.26237017060353,-17.26237017060353,-17.26237017060353,1.63893574083879,-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353",
    "3.825077778895709,-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353",
    "-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353,1.264056473470967,-17.26237017060353,1.95290923181567",
    "-17.26237017060353,-17.26237017060353,-17.26237017060353,3.204814444466126,-17.26237017060353,2.922414390977227,-17.26237017060353,-17.26237017060353,-17.26237017060353,1.882280604637242",
    "-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353",
    "-17.26237017060353,-17.26237017060353,-17.26237017060353,2.348101417408984,-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353,2.087509263465908,1.848686518557616",
    "-17.26237017060353,-17.26237017060353,0.7657355386920179,-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353,2.213903362554662",
    "-17.26237017060353,-17.26237017060353,0.2964765077157156,-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353",
    "-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353,2.138461329552182,-17.26237017060353",
    "-17.26237017060353,-17.26237017060353,2.150582080611352,-17.26237017060353,-17.26237017060353,-17.26237017060353,1.242964571461979,2.07902051742456,2.105021116863338,-17.26237017060353",
    "-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353,-17.26237017060353,1.525892533090062,1.023361652838154,1.547463887562844,-17.26237017060353",
    "-17.26237017060353,-17.26237017060353,1.967952890748702,-17.26237017060353,-17.26237017060353,-17.26237017060353,1.542621052141231,1.585999087662979,-17.26237017060353,-17.26237017060353",
    "-17.26237017060353,-17.26237017060353,1.653477121538708,-17.26237017060353,3.175273209550256,-17.26237017060353,2.032809087343519,0.8022936630309349,-17.26237017060353,2.070656299484464",
    "-17.26237017060353,-17.26237017060353,-17.26237017060353,1.363841502151723,-17.26237017060353,2.487657972269282,2.132613695219045,-17.26237017060353,-17.26237017060353,-17.26237017060353"
)

$nRows = $rowData.Count
$nCols = ($rowData[0] -split ",").Count

$values = New-Object 'object[,]' $nRows, $nCols
for ($r = 0; $r -lt $nRows; $r++) {
    $cells = $rowData[$r] -split ","
    for ($c = 0; $c -lt $nCols; $c++) {
        $values[$r, $c] = [double]$cells[$c]
    }
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2:K21").Value = $values
